$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 659.8333
$ws.Range("I6").Value = 213.11111
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 639.3333299999999
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = -527.3333299999999
$ws.Range("N6").Value = -6224
$ws.Range("H98").Value = 490492.1
$ws.Range("I98").Value = 625712
$ws.Range("K98").Value = 625712
$ws.Range("M98").Value = -624214
$ws.Range("H122").Value = 490492.1
$ws.Range("I122").Value = 625712
$ws.Range("K122").Value = 1877136
$ws.Range("M122").Value = -1874686
$ws.Range("H125").Value = 14015131
$ws.Range("I125").Value = 300
$ws.Range("J125").Value = 16017250
$ws.Range("K125").Value = 2700
$ws.Range("L125").Value = 144155250
$ws.Range("M125").Value = -240
$ws.Range("N125").Value = -144160170
$ws.Range("H129").Value = 1147.5161
$ws.Range("I129").Value = 363.33334
$ws.Range("J129").Value = 1231.5358
$ws.Range("K129").Value = 1090.00002
$ws.Range("L129").Value = 3694.6074
$ws.Range("M129").Value = 3909.99998
$ws.Range("N129").Value = -13694.6074
$ws.Range("H132").Value = 25309.12
$ws.Range("I132").Value = 27197.104
$ws.Range("J132").Value = 765.3333
$ws.Range("K132").Value = 81591.31200000001
$ws.Range("L132").Value = 2295.9999
$ws.Range("M132").Value = -79061.31200000001
$ws.Range("N132").Value = -7355.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 2400
$ws.Range("I13").Value = 1500
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = -1356
$ws.Range("N13").Value = -3288
$ws.Range("H32").Value = 4208.11
$ws.Range("I32").Value = 2773
$ws.Range("J32").Value = 15819.454
$ws.Range("K32").Value = 2773
$ws.Range("L32").Value = 15819.454
$ws.Range("M32").Value = -2486
$ws.Range("N32").Value = -16393.454
$ws.Range("H53").Value = 2000
$ws.Range("I53").Value = 2000
$ws.Range("K53").Value = 2000
$ws.Range("M53").Value = -1318
$ws.Range("H61").Value = 1738.619
$ws.Range("I61").Value = 1263.8286
$ws.Range("K61").Value = 1263.8286
$ws.Range("M61").Value = -1051.8286
$ws.Range("H136").Value = 1738.619
$ws.Range("I136").Value = 1263.8286
$ws.Range("K136").Value = 3791.4858
$ws.Range("M136").Value = -1241.4858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2447.0278
$ws.Range("I134").Value = 1452.6
$ws.Range("J134").Value = 7419.1665
$ws.Range("K134").Value = 4357.799999999999
$ws.Range("L134").Value = 22257.4995
$ws.Range("M134").Value = -1822.799999999999
$ws.Range("N134").Value = -27327.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4675.6333
$ws.Range("I31").Value = 1789.0605
$ws.Range("K31").Value = 1789.0605
$ws.Range("M31").Value = -1494.0605
$ws.Range("H32").Value = 30000
$ws.Range("J32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("N32").Value = -20632
$ws.Range("H34").Value = 4675.6333
$ws.Range("I34").Value = 1789.0605
$ws.Range("K34").Value = 1789.0605
$ws.Range("M34").Value = -1587.0605
$ws.Range("H58").Value = 2694.8
$ws.Range("I58").Value = 1389.2727
$ws.Range("K58").Value = 1389.2727
$ws.Range("M58").Value = -1186.2727
$ws.Range("H136").Value = 2694.8
$ws.Range("I136").Value = 1389.2727
$ws.Range("K136").Value = 4167.8181
$ws.Range("M136").Value = -1617.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 1227.6666
$ws.Range("I28").Value = 725
$ws.Range("J28").Value = 2233
$ws.Range("K28").Value = 2175
$ws.Range("L28").Value = 6699
$ws.Range("M28").Value = -1943
$ws.Range("N28").Value = -7163
$ws.Range("H99").Value = 1554
$ws.Range("I99").Value = 1081
$ws.Range("K99").Value = 3243
$ws.Range("M99").Value = -997
$ws.Range("H113").Value = 548.9286
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 554.8
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1664.4
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6004.4
$ws.Range("H122").Value = 1125.3636
$ws.Range("I122").Value = 326.9
$ws.Range("J122").Value = 1790.75
$ws.Range("K122").Value = 2942.1
$ws.Range("L122").Value = 16116.75
$ws.Range("M122").Value = -492.0999999999999
$ws.Range("N122").Value = -21016.75
$ws.Range("H134").Value = 5411.654
$ws.Range("I134").Value = 1850.25
$ws.Range("J134").Value = 11109.9
$ws.Range("K134").Value = 5550.75
$ws.Range("L134").Value = 33329.7
$ws.Range("M134").Value = -480.75
$ws.Range("N134").Value = -43469.7
$ws.Range("H137").Value = 9186211
$ws.Range("I137").Value = 16668204
$ws.Range("J137").Value = 207819.6
$ws.Range("K137").Value = 50004612
$ws.Range("L137").Value = 623458.8
$ws.Range("M137").Value = -49999512
$ws.Range("N137").Value = -633658.8
$ws.Range("H141").Value = 5682.857
$ws.Range("I141").Value = 8070
$ws.Range("J141").Value = 2500
$ws.Range("K141").Value = 24210
$ws.Range("L141").Value = 7500
$ws.Range("M141").Value = -19030
$ws.Range("N141").Value = -17860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 25250
$ws.Range("I10").Value = 25250
$ws.Range("K10").Value = 25250
$ws.Range("M10").Value = -25081

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 50000
$ws.Range("I39").Value = 50000
$ws.Range("K39").Value = 50000
$ws.Range("M39").Value = -49540
$ws.Range("H61").Value = 1869.7333
$ws.Range("I61").Value = 1626.4546
$ws.Range("J61").Value = 2538.75
$ws.Range("K61").Value = 1626.4546
$ws.Range("L61").Value = 2538.75
$ws.Range("M61").Value = -1424.4546
$ws.Range("N61").Value = -2942.75
$ws.Range("H74").Value = 193464.67
$ws.Range("I74").Value = 275197
$ws.Range("J74").Value = 30000
$ws.Range("K74").Value = 275197
$ws.Range("L74").Value = 30000
$ws.Range("M74").Value = -274199
$ws.Range("N74").Value = -31996
$ws.Range("H77").Value = 193464.67
$ws.Range("I77").Value = 275197
$ws.Range("J77").Value = 30000
$ws.Range("K77").Value = 825591
$ws.Range("L77").Value = 90000
$ws.Range("M77").Value = -820599
$ws.Range("N77").Value = -99984
$ws.Range("H93").Value = 1442.9375
$ws.Range("I93").Value = 1410.5
$ws.Range("J93").Value = 1497
$ws.Range("K93").Value = 1410.5
$ws.Range("L93").Value = 1497
$ws.Range("M93").Value = -162.5
$ws.Range("N93").Value = -3993
$ws.Range("H101").Value = 29900
$ws.Range("J101").Value = 29900
$ws.Range("L101").Value = 29900
$ws.Range("N101").Value = -36390
$ws.Range("H113").Value = 1869.7333
$ws.Range("I113").Value = 1626.4546
$ws.Range("J113").Value = 2538.75
$ws.Range("K113").Value = 1626.4546
$ws.Range("L113").Value = 2538.75
$ws.Range("M113").Value = 543.5454
$ws.Range("N113").Value = -6878.75
$ws.Range("H136").Value = 4322.88
$ws.Range("I136").Value = 2261.6667
$ws.Range("J136").Value = 15144.25
$ws.Range("K136").Value = 6785.000100000001
$ws.Range("L136").Value = 45432.75
$ws.Range("M136").Value = -4235.000100000001
$ws.Range("N136").Value = -50532.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 130
$ws.Range("I8").Value = 130
$ws.Range("K8").Value = 130
$ws.Range("M8").Value = 10
$ws.Range("H10").Value = 2333.6667
$ws.Range("I10").Value = 1001
$ws.Range("K10").Value = 1001
$ws.Range("M10").Value = -832
$ws.Range("H95").Value = 24900
$ws.Range("J95").Value = 24900
$ws.Range("L95").Value = 24900
$ws.Range("N95").Value = -30392
$ws.Range("H103").Value = 517801
$ws.Range("J103").Value = 517801
$ws.Range("L103").Value = 517801
$ws.Range("N103").Value = -520145
$ws.Range("H132").Value = 9120.299999999999
$ws.Range("I132").Value = 10175.5
$ws.Range("J132").Value = 4899.5
$ws.Range("K132").Value = 30526.5
$ws.Range("L132").Value = 14698.5
$ws.Range("M132").Value = -27996.5
$ws.Range("N132").Value = -19758.5
$ws.Range("H136").Value = 1792
$ws.Range("I136").Value = 1069.6316
$ws.Range("J136").Value = 4079.5
$ws.Range("K136").Value = 3208.8948
$ws.Range("L136").Value = 12238.5
$ws.Range("M136").Value = -658.8948
$ws.Range("N136").Value = -17338.5
